$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for column E
$ws.Range("E1").Value = "strength (RMS)"

# Update data rows 2-19 for columns B, D, E (C is unchanged)
$data = @(
    @{Row=2;  B=16.4; D=14.4; E=50.4},
    @{Row=3;  B=15.6; D=14.8; E=44.4},
    @{Row=4;  B=15.6; D=15.2; E=40.4},
    @{Row=5;  B=16.8; D=14.8; E=44.4},
    @{Row=6;  B=18;   D=12.4; E=42},
    @{Row=7;  B=15.6; D=15.6; E=41.4},
    @{Row=8;  B=16;   D=14;   E=41.67},
    @{Row=9;  B=15.2; D=15.2; E=42.8},
    @{Row=10; B=16.8; D=14;   E=40.4},
    @{Row=11; B=16;   D=14.4; E=39.6},
    @{Row=12; B=16.4; D=14;   E=36.6},
    @{Row=13; B=16;   D=14;   E=44.8},
    @{Row=14; B=15.6; D=15.2; E=42.2},
    @{Row=15; B=15.6; D=14.8; E=46.4},
    @{Row=16; B=21.6; D=8.800000000000001; E=35},
    @{Row=17; B=19.2; D=12.8; E=43.4},
    @{Row=18; B=15.5; D=15;   E=48},
    @{Row=19; B=15.6; D=14.8; E=47.4}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
